# Weekly price-sheet update: two new "Apio" (celery) price records are
# inserted into the middle of the data block (not merely appended), which
# pushes the existing rows below them down. This mirrors the upstream
# source feed, which keeps rows ordered by something other than simple
# append order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues([int]$Row, [object[]]$Values) {
    for ($i = 0; $i -lt $Values.Count; $i++) {
        $ws.Cells.Item($Row, $i + 1).Value2 = $Values[$i]
    }
}

# --- Insertion #1: one new row, becomes row 263 -------------------------
$ws.Rows.Item(263).Insert()

Set-RowValues 263 @(
    11,
    "Vega Monumental Concepción",
    "Bíobío",
    44748,
    8,
    100112017,
    "Apio",
    "Americana (o)",
    "Primera",
    150,
    7000,
    7500,
    7167,
    "`$/docena de matas",
    "Región de Coquimbo",
    1194,
    6,
    "Hortaliza"
)

# --- Insertion #2: two new rows, become rows 280 and 281 ----------------
$ws.Range("A280:A281").EntireRow.Insert()

Set-RowValues 280 @(
    11,
    "Vega Monumental Concepción",
    "Bíobío",
    44747,
    8,
    100112017,
    "Apio",
    "Americana (o)",
    "Primera",
    270,
    8000,
    8500,
    8222,
    "`$/docena de matas",
    "Región de Coquimbo",
    1370,
    6,
    "Hortaliza"
)

Set-RowValues 281 @(
    11,
    "Vega Monumental Concepción",
    "Bíobío",
    44747,
    8,
    100112017,
    "Apio",
    "Americana (o)",
    "Segunda",
    130,
    7000,
    7000,
    7000,
    "`$/docena de matas",
    "Región de Coquimbo",
    1167,
    6,
    "Hortaliza"
)
